# Fix data that was incorrectly pulled from other tickers' files.
# For every data row (2-43) on Sheet1:
#   - column I (fixed_ticker)      -> "TEMN SW"
#   - column AL (Original Currency) -> "CHF"   (unchanged value, but re-affirmed)
#   - columns D/E/F/G/H (open/close/high/low price, shares_outstanding)
#     get replaced with the correct TEMENOS AG REG (TEMN SW) figures.
# Rows 2-42 all share one constant set of corrected values; row 43 (the most
# recent date) has its own distinct corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$openPrice  = 77.72841848241376
$closePrice = 72.10655346679688
$highPrice  = 79.07277608898339
$lowPrice   = 68.68454847198167
$shares     = 68984148

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 4).Value = $openPrice
    $ws.Cells.Item($r, 5).Value = $closePrice
    $ws.Cells.Item($r, 6).Value = $highPrice
    $ws.Cells.Item($r, 7).Value = $lowPrice
    $ws.Cells.Item($r, 8).Value = $shares
    $ws.Cells.Item($r, 9).Value = "TEMN SW"
    $ws.Cells.Item($r, 38).Value = "CHF"
}

# Row 43 has its own distinct corrected values.
$ws.Cells.Item(43, 4).Value = 70.824
$ws.Cells.Item(43, 5).Value = 91.29120190429687
$ws.Cells.Item(43, 6).Value = 93.16320190429687
$ws.Cells.Item(43, 7).Value = 69.264
$ws.Cells.Item(43, 8).Value = 68984148
$ws.Cells.Item(43, 9).Value = "TEMN SW"
$ws.Cells.Item(43, 38).Value = "CHF"
